# Sprint 2 Backlog (day 4) update
# Fills in the "DAY 5" (column J) remaining-hours for every user story row,
# bumps a couple of "DAY 4" (column I) values that were re-estimated, and
# refreshes the totals row (row 48) that feed the burndown chart. Also
# reselects the cell the author ended the session on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Column J ("DAY 5") remaining-hours entries for the backlog rows whose
# cell formatting does not change (style stays exactly as it was).
# -----------------------------------------------------------------
$dayFiveValues = [ordered]@{
    3  = 0
    5  = 0
    7  = 0
    9  = 2
    10 = 0
    11 = 2
    12 = 0
    13 = 4
    14 = 0
    15 = 3
    16 = 1
    17 = 3
    18 = 2
    19 = 1
    20 = 4
    21 = 0
    22 = 3
    23 = 1
    24 = 3
    25 = 0
    26 = 2
    27 = 1
    29 = 4
    30 = 0
    31 = 3
    32 = 1
    33 = 4
    34 = 3
    35 = 1
    36 = 4
    37 = 3
    38 = 1
    39 = 3
    40 = 2
    41 = 1
    43 = 2
    44 = 0
    46 = 1
}

foreach ($row in $dayFiveValues.Keys) {
    $ws.Range("J$row").Value = $dayFiveValues[$row]
}

# -----------------------------------------------------------------
# Rows 42 and 45 are the two "story" summary rows where column J had not
# been touched before, so Excel re-used the "DAY 4" cell's (column I)
# left-aligned style when the value was typed in, instead of the
# previously-unused centred style.
# -----------------------------------------------------------------
$ws.Range("J42").HorizontalAlignment = -4131   # xlLeft, matches I42's style
$ws.Range("J42").Value = 2

$ws.Range("I45").Value = 3
$ws.Range("J45").HorizontalAlignment = -4131   # xlLeft, matches I45's style
$ws.Range("J45").Value = 1.5

# Row 47: DAY 4 re-estimated from 0.5 to 1, and DAY 5 now has its first entry.
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 0.5

# -----------------------------------------------------------------
# Row 48 holds the manually maintained totals (DAY 4 / DAY 5 columns)
# that the burndown chart reads from.
# -----------------------------------------------------------------
$ws.Range("I48").Value = 44
$ws.Range("J48").Value = 34.5

# -----------------------------------------------------------------
# Restore the author's final selection on the sheet.
# -----------------------------------------------------------------
$ws.Range("B50:V50").Select()
